$d = $word.ActiveDocument

$pairs = @(
    @("柴茶：生活的香料", "印度奶茶：生活的香料"),
    @("柴茶：杯子里的味道世界", "印度奶茶：杯子里的味道世界"),
    @("柴茶：发现印度的魔力", "印度奶茶：感受印度的魔力"),
    @("柴茶：健康与快乐的完美融合", "印度奶茶：健康与快乐的完美融合"),
    @("柴茶：不仅仅是茶，一种生活方式", "印度奶茶：不仅仅是茶，更是一种生活方式"),
    @("柴茶：所有季节和原因的饮料", "印度奶茶：适合所有季节的饮品，品尝无需原因"),
    @("柴茶：你的感官的终极放纵", "印度奶茶：感官的终极放纵"),
    @("柴茶：从日常的甜蜜逃跑", "印度奶茶：感受日常流露的甜蜜"),
    @("柴茶：分享温暖，分享爱", "印度奶茶：分享温暖，分享爱"),
    @("柴茶：把自己当作一些特别的东西", "印度奶茶：款待自己一些特别的东西")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
